# Trait Matrix workbook edit
# - Adds a new "Trait codings" sheet (after "Comments") that provides a
#   human-readable lookup of trait name -> long name -> coded value labels.
# - Makes the new sheet the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last sheet (after "Comments") -----------
$commentsSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $commentsSheet)
$ws.Name = "Trait codings"

# --- Header row -------------------------------------------------------------
$ws.Range("A1").Value = "Trait name"
$ws.Range("B1").Value = "Long name"
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 1
$ws.Range("E1").Value = 2

# --- Trait code lookup rows --------------------------------------------------
$ws.Range("A2").Value = "Size"
$ws.Range("B2").Value = "Valve Size"
$ws.Range("C2").Value = "Small (< 0.50mm x 0.25mm)"
$ws.Range("D2").Value = "Medium (0.50mm-1.00mm x 0.25mm " + [char]0x2013 + " 0.50mm)"
$ws.Range("E2").Value = "Large (> 1.00mm x 0.50mm)"

$ws.Range("A3").Value = "Calcification"
$ws.Range("B3").Value = "Valve Calcification"
$ws.Range("C3").Value = "Weak"
$ws.Range("D3").Value = "Average"
$ws.Range("E3").Value = "Thin"

$ws.Range("A4").Value = "Shape"
$ws.Range("B4").Value = "Valve Shape (dorsal view: anterior-posterior)"
$ws.Range("C4").Value = "Ovate or Subovate"
$ws.Range("D4").Value = "Rectangular or Subrectangular"

$ws.Range("A5").Value = "Carapace Texture"
$ws.Range("B5").Value = "Presence of Surface Reticulation"
$ws.Range("C5").Value = "Absent"
$ws.Range("D5").Value = "Present"

$ws.Range("A6").Value = "Ventral Margin"
$ws.Range("B6").Value = "Ventral Margin Shape"
$ws.Range("C6").Value = "Convex"
$ws.Range("D6").Value = "Straight"
$ws.Range("E6").Value = "Concave"

$ws.Range("A7").Value = "Dorsal Margin"
$ws.Range("B7").Value = "Dorsal Margin Shape"
$ws.Range("C7").Value = "Convex"
$ws.Range("D7").Value = "Straight"
$ws.Range("E7").Value = "Concave"

$ws.Range("A8").Value = "Posterior Margin"
$ws.Range("B8").Value = "Posterior Margin Shape"
$ws.Range("C8").Value = "Convex"
$ws.Range("D8").Value = "Straight"
$ws.Range("E8").Value = "Concave"

$ws.Range("A9").Value = "Anterior Margin"
$ws.Range("B9").Value = "Anterior Margin Shape"
$ws.Range("C9").Value = "Convex"
$ws.Range("D9").Value = "Straight"
$ws.Range("E9").Value = "Concave"

$ws.Range("A10").Value = "Left overlap"
$ws.Range("B10").Value = "Right/Left Valve Size Ratio"
$ws.Range("C10").Value = "Left Valve Larger"
$ws.Range("D10").Value = "Equally Sized Valves"
$ws.Range("E10").Value = "Right Valve Larger"
$ws.Range("A10:B10").Font.Bold = $true

$ws.Range("A11").Value = "Spines"
$ws.Range("B11").Value = "Presence of Spines"
$ws.Range("C11").Value = "Absent"
$ws.Range("D11").Value = "Present"

$ws.Range("A12").Value = "Carapace pits"
$ws.Range("B12").Value = "Presence of Carapace Pits"
$ws.Range("C12").Value = "Absent"
$ws.Range("D12").Value = "Present"
$ws.Range("A12:B12").Font.Bold = $true

$ws.Range("A13").Value = "Opaque Areas"
$ws.Range("B13").Value = "Presence of Opaque Patches on Valve"
$ws.Range("C13").Value = "Absent"
$ws.Range("D13").Value = "Present"
$ws.Range("A13:B13").Font.Bold = $true

$ws.Range("A14").Value = "Denticulations"
$ws.Range("B14").Value = "Presence of Denticulations"
$ws.Range("C14").Value = "Absent"
$ws.Range("D14").Value = "Present"

$ws.Range("A15").Value = "ala"
$ws.Range("B15").Value = "Presence of Alae"
$ws.Range("C15").Value = "Absent"
$ws.Range("D15").Value = "Present"

$ws.Range("A16").Value = "nodes"
$ws.Range("B16").Value = "Presence of Nodes"
$ws.Range("C16").Value = "Absent"
$ws.Range("D16").Value = "Present"

$ws.Range("A17").Value = "caudal process"
$ws.Range("B17").Value = "Presence of Caudal Process"
$ws.Range("C17").Value = "Absent"
$ws.Range("D17").Value = "Present"

$ws.Range("A18").Value = "sulcus"
$ws.Range("B18").Value = "Presence of Sulcus"
$ws.Range("C18").Value = "Absent"
$ws.Range("D18").Value = "Present"

$ws.Range("A19").Value = "eye tubercule"
$ws.Range("B19").Value = "Presence of Eye Tubercules"
$ws.Range("C19").Value = "Absent"
$ws.Range("D19").Value = "Present"

# --- Cosmetic column widths on the new sheet --------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.9392712550607
$ws.Columns.Item(2).ColumnWidth = 27.7004048582996
$ws.Columns.Item(3).ColumnWidth = 27.7004048582996
$ws.Columns.Item(4).ColumnWidth = 45.6113360323887
$ws.Columns.Item(5).ColumnWidth = 27.7651821862348

# --- Select the first cell and make this new sheet the active tab ----------
$ws.Range("A12").Select()
$ws.Activate()
